$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells we touch to stay text-typed so
# numeric-looking strings (e.g. "1.00", "61.80") keep their exact
# textual representation instead of being coerced to a number.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.414.87'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.539.84'
$ws.Range('E3').Value = '  -1.79%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '195.52'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '582.53'
$ws.Range('E6').Value = '  -3.65%  '
$ws.Range('E7').Value = '  -2.48%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('E10').Value = '  -2.82%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '51.73'
$ws.Range('E11').Value = '  -3.82%  '
$ws.Range('E12').Value = '  -5.74%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '9.23'
$ws.Range('E13').Value = '  -3.85%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.097.16'
$ws.Range('E14').Value = '  -1.88%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '665.98'
$ws.Range('E15').Value = '  +11.90%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '69.536.97'
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '12.53'
$ws.Range('E17').Value = '  -4.30%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.524.16'
$ws.Range('E18').Value = '  -2.23%  '
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.44'
$ws.Range('E20').Value = '  -3.46%  '
$ws.Range('E21').Value = '  -3.24%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '18.17'
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.31'
$ws.Range('E23').Value = '  +2.58%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '104.27'
$ws.Range('E24').Value = '  +2.46%  '
$ws.Range('E25').Value = '  -5.20%  '
$ws.Range('E26').Value = '  -4.00%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.17'
$ws.Range('E27').Value = '  -5.26%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.58'
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '33.11'
$ws.Range('E29').Value = '  -1.95%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.35'
$ws.Range('E30').Value = '  -7.35%  '
$ws.Range('E31').Value = '  -5.57%  '
$ws.Range('E32').Value = '  -4.29%  '
$ws.Range('E33').Value = '  -5.24%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '61.80'
$ws.Range('E34').Value = '  -2.36%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.789.12'
$ws.Range('E35').Value = '  -3.03%  '
$ws.Range('E36').Value = '  -7.73%  '
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.69'
$ws.Range('E38').Value = '  +4.40%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '500.29'
$ws.Range('E39').Value = '  -3.31%  '
$ws.Range('E40').Value = '  -6.48%  '
$ws.Range('E41').Value = '  -5.12%  '
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '34.65'
$ws.Range('E43').Value = '  -5.98%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0448'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.36'
$ws.Range('E45').Value = '  -2.24%  '
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('E47').Value = '  -2.79%  '
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.31'
$ws.Range('E49').Value = '  -3.65%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.76'
$ws.Range('E50').Value = '  +19.35%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.70'
$ws.Range('E51').Value = '  +63.28%  '
